# Plan_v1.xlsx update:
#  - Several checklist requirements that were previously "not started" / "in
#    progress" are now completed (Status column D -> 100%), which also
#    switches their cell style from the "Bad"/"Neutral" (red/yellow) look to
#    the "Good" (green) look already used by the other completed rows.
#  - The current selection/view is moved further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell that already carries the "Good" (green, percent, centered)
# status style - reuse it so no new style entries are introduced.
$goodTemplate = $ws.Range("D3")
$goodTemplate.Copy()

# Rows whose requirement got finished (Status 0 / 0.6 -> 1 = 100%).
$doneCells = @("D8", "D10", "D12", "D17", "D18")
foreach ($cellAddr in $doneCells) {
    $cell = $ws.Range($cellAddr)
    $cell.PasteSpecial(-4122) # xlPasteFormats - copy the "Good" look
    $cell.Value = 1
}

$excel.CutCopyMode = 0

# Move the active selection further down the list.
$ws.Activate()
$ws.Range("F10").Select()
